function New-BodyXml($innerXml) {
    return '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:body>' + $innerXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

$d = $word.ActiveDocument

# Apply edits from the END of the document toward the START, so
# earlier character offsets are not invalidated by later structural edits.

# --- A. PATO hyperlink paragraph (list item 56) ---
$p = $d.Paragraphs.Item(56)
$rng = $d.Range($p.Range.Start, $p.Range.End)
$rng.InsertXML((New-BodyXml('<w:p w:rsidR="00E25B91" w:rsidRDefault="0086107B"><w:pPr><w:widowControl w:val="0"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:ind w:hanging="360"/><w:contextualSpacing/><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:val="pt-BR"/></w:rPr></w:pPr><w:hyperlink r:id="rId47"><w:proofErr w:type="gramStart"/><w:r w:rsidR="005B1048"><w:rPr><w:color w:val="1155CC"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:u w:val="single"/><w:lang w:val="pt-BR"/></w:rPr><w:t>BAO2.</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r w:rsidR="005B1048"><w:rPr><w:color w:val="1155CC"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:u w:val="single"/><w:lang w:val="pt-BR"/></w:rPr><w:t>3</w:t></w:r><w:r w:rsidR="005B1048"><w:rPr><w:color w:val="1155CC"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:u w:val="single"/><w:lang w:val="pt-BR"/></w:rPr><w:t>0_PATO_import.owl</w:t></w:r></w:hyperlink></w:p>')))

# --- B. NCBITaxon hyperlink paragraph (list item 55) ---
$p = $d.Paragraphs.Item(55)
$rng = $d.Range($p.Range.Start, $p.Range.End)
$rng.InsertXML((New-BodyXml('<w:p w:rsidR="00E25B91" w:rsidRDefault="0086107B"><w:pPr><w:widowControl w:val="0"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:ind w:hanging="360"/><w:contextualSpacing/><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:val="pt-BR"/></w:rPr></w:pPr><w:hyperlink r:id="rId46"><w:proofErr w:type="gramStart"/><w:r w:rsidR="005B1048"><w:rPr><w:color w:val="1155CC"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:u w:val="single"/><w:lang w:val="pt-BR"/></w:rPr><w:t>BAO2.</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r w:rsidR="005B1048"><w:rPr><w:color w:val="1155CC"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:u w:val="single"/><w:lang w:val="pt-BR"/></w:rPr><w:t>3</w:t></w:r><w:r w:rsidR="005B1048"><w:rPr><w:color w:val="1155CC"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:u w:val="single"/><w:lang w:val="pt-BR"/></w:rPr><w:t>_NCBITaxon_import.owl</w:t></w:r></w:hyperlink></w:p>')))

# --- C. DOID hyperlink paragraph (list item 53) ---
$p = $d.Paragraphs.Item(53)
$rng = $d.Range($p.Range.Start, $p.Range.End)
$rng.InsertXML((New-BodyXml('<w:p w:rsidR="00E25B91" w:rsidRDefault="0086107B"><w:pPr><w:widowControl w:val="0"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:ind w:hanging="360"/><w:contextualSpacing/><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:val="pt-BR"/></w:rPr></w:pPr><w:hyperlink r:id="rId44"><w:proofErr w:type="gramStart"/><w:r w:rsidR="005B1048"><w:rPr><w:color w:val="1155CC"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:u w:val="single"/><w:lang w:val="pt-BR"/></w:rPr><w:t>BAO2.</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r w:rsidR="005B1048"><w:rPr><w:color w:val="1155CC"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:u w:val="single"/><w:lang w:val="pt-BR"/></w:rPr><w:t>3</w:t></w:r><w:r w:rsidR="005B1048"><w:rPr><w:color w:val="1155CC"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:u w:val="single"/><w:lang w:val="pt-BR"/></w:rPr><w:t>_DOID_import.owl</w:t></w:r></w:hyperlink></w:p>')))

# --- D. CLO hyperlink paragraph (list item 52) ---
$p = $d.Paragraphs.Item(52)
$rng = $d.Range($p.Range.Start, $p.Range.End)
$rng.InsertXML((New-BodyXml('<w:p w:rsidR="00E25B91" w:rsidRDefault="0086107B"><w:pPr><w:widowControl w:val="0"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:ind w:hanging="360"/><w:contextualSpacing/><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:val="pt-BR"/></w:rPr></w:pPr><w:hyperlink r:id="rId43"><w:proofErr w:type="gramStart"/><w:r w:rsidR="005B1048"><w:rPr><w:color w:val="1155CC"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:u w:val="single"/><w:lang w:val="pt-BR"/></w:rPr><w:t>BAO2.</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r w:rsidR="005B1048"><w:rPr><w:color w:val="1155CC"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:u w:val="single"/><w:lang w:val="pt-BR"/></w:rPr><w:t>3</w:t></w:r><w:r w:rsidR="005B1048"><w:rPr><w:color w:val="1155CC"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:u w:val="single"/><w:lang w:val="pt-BR"/></w:rPr><w:t>_CLO_import.owl</w:t></w:r></w:hyperlink></w:p>')))

# --- E. Release date run: "27, 2017" -> "30" + ", 2017" ---
$p3 = $d.Paragraphs.Item(3)
$full = $p3.Range.Text
$relStart = $p3.Range.Start + $full.IndexOf("27, 2017")
$relEnd = $relStart + "27, 2017".Length
$rng = $d.Range($relStart, $relEnd)
$rng.InsertXML((New-BodyXml('<w:p><w:r><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>30</w:t></w:r><w:r><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>, 2017</w:t></w:r></w:p>')))

# --- F. Region: About BAO body .. License/Disclaimer heading ---
$p8 = $d.Paragraphs.Item(8)
$p14 = $d.Paragraphs.Item(14)
$rng = $d.Range($p8.Range.Start, $p14.Range.End)
$regionXml = @'
<w:p w:rsidR="00E25B91" w:rsidRDefault="005B1048">
  <w:pPr>
    <w:widowControl w:val="0"/>
    <w:rPr>
      <w:color w:val="333333"/>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:color w:val="333333"/>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
      <w:highlight w:val="white"/>
    </w:rPr>
    <w:t xml:space="preserve">The </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:color w:val="333333"/>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
      <w:highlight w:val="white"/>
    </w:rPr>
    <w:t>BioAssay</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:color w:val="333333"/>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
      <w:highlight w:val="white"/>
    </w:rPr>
    <w:t xml:space="preserve"> Ontology (BAO) has been developed to formally describe biological screening assays and their results including high-throughput screening (HTS) data; specifically in the domain of small molecule drug and probe development. BAO enables categorization of assays and results by based on several concepts that are important to interpret and analyze screening data with the goal to infer the mechanism of action of small molecules based on the known aggregate screening results from many assays.</w:t>
  </w:r>
</w:p>
<w:p w:rsidR="00E25B91" w:rsidRDefault="00E25B91">
  <w:pPr>
    <w:widowControl w:val="0"/>
    <w:rPr>
      <w:color w:val="333333"/>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p>
  <w:pPr>
    <w:widowControl w:val="0"/>
    <w:rPr>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
    </w:rPr>
    <w:t xml:space="preserve">BAO 2.3.1 is based on the need after the annotation effort performed by Collaborative Drug Discovery (CDD) using their annotation tool </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
    </w:rPr>
    <w:t>BioAssay</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
    </w:rPr>
    <w:t xml:space="preserve"> Express (</w:t>
  </w:r>
  <w:proofErr w:type="gramStart"/>
  <w:r>
    <w:rPr>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
    </w:rPr>
    <w:t>BAE :</w:t>
  </w:r>
  <w:proofErr w:type="gramEnd"/>
  <w:r>
    <w:rPr>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
    </w:rPr>
    <w:t xml:space="preserve"> http://www.bioassayexpress.com/). CDD annotated 3500 assays using BAE.</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:widowControl w:val="0"/>
    <w:rPr>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
    </w:rPr>
    <w:t xml:space="preserve">Based on the annotations, absent term reports have been generated by BAE. The absence reports were then communicated with the BAO group and processed systematically in a semi-automated way for the </w:t>
  </w:r>
  <w:bookmarkStart w:id="3" w:name="_GoBack"/>
  <w:bookmarkEnd w:id="3"/>
  <w:r>
    <w:rPr>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
    </w:rPr>
    <w:t>ontology additions. The generated OWL vocabulary files with new terms are added to the different vocabulary files based on BAO’s modular architecture. New terms from BAE were added to BAO 2.2.2</w:t>
  </w:r>
</w:p>
<w:p w:rsidR="00E25B91" w:rsidRDefault="005B1048">
  <w:pPr>
    <w:pStyle w:val="Heading3"/>
    <w:keepNext w:val="0"/>
    <w:keepLines w:val="0"/>
    <w:widowControl w:val="0"/>
    <w:contextualSpacing w:val="0"/>
  </w:pPr>
  <w:bookmarkStart w:id="4" w:name="_8wnlhsiqhnm3" w:colFirst="0" w:colLast="0"/>
  <w:bookmarkEnd w:id="4"/>
  <w:r>
    <w:t>Changes</w:t>
  </w:r>
</w:p>
<w:p w:rsidR="00E25B91" w:rsidRDefault="0086107B">
  <w:pPr>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="1"/>
    </w:numPr>
    <w:ind w:hanging="360"/>
    <w:contextualSpacing/>
  </w:pPr>
  <w:bookmarkStart w:id="5" w:name="_ofe3vws0y5rj" w:colFirst="0" w:colLast="0"/>
  <w:bookmarkEnd w:id="5"/>
  <w:r>
    <w:t>01/30/2015</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve">: A total number of 224 new terms added </w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve">to BAO 2.2.2 </w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve">via annotations generated with </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:t>BioAssay</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:t xml:space="preserve"> Express (BAE). </w:t>
  </w:r>
</w:p>
<w:p w:rsidR="0086107B" w:rsidRDefault="0086107B" w:rsidP="0086107B">
  <w:pPr>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="1"/>
    </w:numPr>
    <w:ind w:hanging="360"/>
    <w:contextualSpacing/>
  </w:pPr>
  <w:r>
    <w:t>Vocabulary and Modules changed are as follows:</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="1"/>
    </w:numPr>
    <w:contextualSpacing/>
  </w:pPr>
  <w:r>
    <w:t xml:space="preserve">BAO Core: </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:t>bao_vocabulary_assay</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:t xml:space="preserve">, </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:t>assaykit</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:t xml:space="preserve">, detection, format, instrument, method, result, and </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:t>screenedentity</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:t>.</w:t>
  </w:r>
</w:p>
<w:p w:rsidR="0086107B" w:rsidRDefault="0086107B" w:rsidP="0086107B">
  <w:pPr>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="1"/>
    </w:numPr>
    <w:contextualSpacing/>
    <w:rPr>
      <w:lang w:val="pt-BR"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:lang w:val="pt-BR"/>
    </w:rPr>
    <w:t xml:space="preserve">BAO </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:lang w:val="pt-BR"/>
    </w:rPr>
    <w:t>External</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:lang w:val="pt-BR"/>
    </w:rPr>
    <w:t xml:space="preserve">: </w:t>
  </w:r>
  <w:proofErr w:type="gramStart"/>
  <w:r>
    <w:rPr>
      <w:lang w:val="pt-BR"/>
    </w:rPr>
    <w:t>BAO2.</w:t>
  </w:r>
  <w:proofErr w:type="gramEnd"/>
  <w:r>
    <w:rPr>
      <w:lang w:val="pt-BR"/>
    </w:rPr>
    <w:t xml:space="preserve">0_DOID_import.owl, BAO2.0_CLO_import.owl, BAO2.0_NCBITaxon_import.owl </w:t>
  </w:r>
</w:p>
<w:p w:rsidR="00E25B91" w:rsidRDefault="005B1048">
  <w:pPr>
    <w:pStyle w:val="Heading3"/>
    <w:keepNext w:val="0"/>
    <w:keepLines w:val="0"/>
    <w:widowControl w:val="0"/>
    <w:contextualSpacing w:val="0"/>
  </w:pPr>
  <w:r>
    <w:t>License/Disclaimer</w:t>
  </w:r>
</w:p>

'@
$rng.InsertXML((New-BodyXml($regionXml)))

# --- G. Title: "2.2" -> "3.1" in the last run of paragraph 1 ---
$p1 = $d.Paragraphs.Item(1)
$full1 = $p1.Range.Text
$tStart = $p1.Range.Start + $full1.Length - 4
$tEnd = $p1.Range.Start + $full1.Length - 1
$rng = $d.Range($tStart, $tEnd)
$rng.InsertXML((New-BodyXml('<w:p><w:r><w:t>3.1</w:t></w:r></w:p>')))

Write-Host "All edits applied"
